$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Price (column D) cells: force text format to avoid numeric auto-conversion ---
$dCells = @("D2","D3","D5","D7","D8","D9","D10","D11","D12","D13","D14","D15","D16","D17","D18","D19","D20","D21","D22","D24","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D46","D47","D48","D49","D50","D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "30.643.72"
$ws.Range("D3").Value = "1.958.62"
$ws.Range("D5").Value = "249.96"
$ws.Range("D7").Value = "0.4830"
$ws.Range("D8").Value = "0.2941"
$ws.Range("D9").Value = "0.06784"
$ws.Range("D10").Value = "110.54"
$ws.Range("D11").Value = "19.48"
$ws.Range("D12").Value = "1.969.27"
$ws.Range("D13").Value = "0.07727"
$ws.Range("D14").Value = "5.470"
$ws.Range("D15").Value = "0.6888"
$ws.Range("D16").Value = "292.49"
$ws.Range("D17").Value = "30.654.93"
$ws.Range("D18").Value = "13.25"
$ws.Range("D19").Value = "2.224.09"
$ws.Range("D20").Value = "5.639"
$ws.Range("D21").Value = "0.000007671"
$ws.Range("D22").Value = "1.001"
$ws.Range("D24").Value = "6.599"
$ws.Range("D25").Value = "9.927"
$ws.Range("D26").Value = "171.25"
$ws.Range("D27").Value = "20.14"
$ws.Range("D28").Value = "2.194"
$ws.Range("D29").Value = "0.1066"
$ws.Range("D30").Value = "1.442"
$ws.Range("D31").Value = "4.706"
$ws.Range("D32").Value = "4.431"
$ws.Range("D33").Value = "0.05098"
$ws.Range("D34").Value = "0.7757"
$ws.Range("D35").Value = "1.173"
$ws.Range("D36").Value = "0.02057"
$ws.Range("D37").Value = "2.733"
$ws.Range("D38").Value = "2.714"
$ws.Range("D39").Value = "2.076"
$ws.Range("D40").Value = "6.265"
$ws.Range("D41").Value = "109.65"
$ws.Range("D42").Value = "0.4467"
$ws.Range("D43").Value = "0.8734"
$ws.Range("D44").Value = "69.95"
$ws.Range("D46").Value = "7.376"
$ws.Range("D47").Value = "0.1279"
$ws.Range("D48").Value = "9.310"
$ws.Range("D49").Value = "35.97"
$ws.Range("D50").Value = "47.59"
$ws.Range("D51").Value = "0.4083"

foreach ($addr in $dCells) {
    $ws.Range($addr).Style = "Normal"
}

# --- Update Volume(1h) (column E) cells ---
$ws.Range("E2").Value = "  +0.43%  "
$ws.Range("E3").Value = "  +2.20%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("E5").Value = "  +1.98%  "
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("E7").Value = "  +0.55%  "
$ws.Range("E8").Value = "  +1.81%  "
$ws.Range("E9").Value = "  +0.84%  "
$ws.Range("E10").Value = "  +0.32%  "
$ws.Range("E11").Value = "  +0.84%  "
$ws.Range("E12").Value = "  +2.89%  "
$ws.Range("E13").Value = "  +2.19%  "
$ws.Range("E14").Value = "  +4.42%  "
$ws.Range("E15").Value = "  +3.16%  "
$ws.Range("E16").Value = "  -3.24%  "
$ws.Range("E17").Value = "  +0.50%  "
$ws.Range("E18").Value = "  +2.08%  "
$ws.Range("E19").Value = "  +2.93%  "
$ws.Range("E20").Value = "  +2.94%  "
$ws.Range("E21").Value = "  +1.42%  "
$ws.Range("E22").Value = "  +0.16%  "
$ws.Range("E23").Value = "  +0.20%  "
$ws.Range("E24").Value = "  +2.66%  "
$ws.Range("E25").Value = "  +4.65%  "
$ws.Range("E26").Value = "  +4.24%  "
$ws.Range("E27").Value = "  -1.27%  "
$ws.Range("E28").Value = "  +4.34%  "
$ws.Range("E29").Value = "  -0.89%  "
$ws.Range("E30").Value = "  +3.33%  "
$ws.Range("E31").Value = "  +17.07%  "
$ws.Range("E32").Value = "  +6.43%  "
$ws.Range("E33").Value = "  +2.38%  "
$ws.Range("E34").Value = "  +5.39%  "
$ws.Range("E35").Value = "  +3.12%  "
$ws.Range("E36").Value = "  +0.30%  "
$ws.Range("E37").Value = "  +0.19%  "
$ws.Range("E38").Value = "  +1.51%  "
$ws.Range("E39").Value = "  +2.95%  "
$ws.Range("E40").Value = "  +6.06%  "
$ws.Range("E41").Value = "  -1.14%  "
$ws.Range("E42").Value = "  +0.96%  "
$ws.Range("E43").Value = "  +0.89%  "
$ws.Range("E44").Value = "  +2.05%  "
$ws.Range("E45").Value = "  +0.18%  "
$ws.Range("E46").Value = "  +1.16%  "
$ws.Range("E47").Value = "  +3.72%  "
$ws.Range("E48").Value = "  +0.65%  "
$ws.Range("E49").Value = "  +3.39%  "
$ws.Range("E50").Value = "  -4.25%  "
$ws.Range("E51").Value = "  +1.82%  "
